# Update "想去人数" (number of people interested) values across the four
# worksheets of the 上海-漫展信息 workbook, per the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 1852
$ws1.Range("F9").Value  = 39
$ws1.Range("F12").Value = 1137
$ws1.Range("F13").Value = 1590
$ws1.Range("F14").Value = 818
$ws1.Range("F15").Value = 1737
$ws1.Range("F19").Value = 192
$ws1.Range("F22").Value = 601
$ws1.Range("F23").Value = 2520
$ws1.Range("F25").Value = 296
$ws1.Range("F31").Value = 4285

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 4165
$ws2.Range("F19").Value = 13
$ws2.Range("F20").Value = 13
$ws2.Range("F31").Value = 4

# --- Sheet 3: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1342
$ws3.Range("F5").Value = 1736
$ws3.Range("F7").Value = 295

# --- Sheet 4: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1342
$ws4.Range("F3").Value  = 1736
$ws4.Range("F5").Value  = 295
$ws4.Range("F12").Value = 1852
$ws4.Range("F15").Value = 39
$ws4.Range("F20").Value = 1137
$ws4.Range("F21").Value = 1590
$ws4.Range("F23").Value = 818
$ws4.Range("F24").Value = 1737
$ws4.Range("F28").Value = 192
$ws4.Range("F31").Value = 13
$ws4.Range("F34").Value = 601
$ws4.Range("F37").Value = 2520
$ws4.Range("F48").Value = 4285
